$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state for rows 2-19 (Player, Position, Team)
$data = @(
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
